$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Designer_Profit  (KPI update for Designer)
# ---------------------------------------------------------------------------
$wsDesigner = $wb.Worksheets.Item("Designer_Profit")
$wsDesigner.Range("B2").Value  = 23.95
$wsDesigner.Range("B3").Value  = 2748.07
$wsDesigner.Range("B4").Value  = 1113.05
$wsDesigner.Range("B5").Value  = 25.64
$wsDesigner.Range("B6").Value  = 67.92
$wsDesigner.Range("B7").Value  = -3.41
$wsDesigner.Range("B9").Value  = 884.83
$wsDesigner.Range("B10").Value = 18.23
$wsDesigner.Range("B11").Value = 247.97
$wsDesigner.Range("B12").Value = 20.69
$wsDesigner.Range("B13").Value = 81.2
$wsDesigner.Range("B14").Value = -1.6

# ---------------------------------------------------------------------------
# Sheet: RD_Profit  (KPI update for R&D; rows 11/12 IDs swap -> TN, JN)
# ---------------------------------------------------------------------------
$wsRD = $wb.Worksheets.Item("RD_Profit")
$wsRD.Range("B2").Value  = 559.9
$wsRD.Range("B3").Value  = 39.37
$wsRD.Range("B5").Value  = 1083.41
$wsRD.Range("B6").Value  = 488.74
$wsRD.Range("B7").Value  = 130.84
$wsRD.Range("B8").Value  = 53.57
$wsRD.Range("B9").Value  = 162.54
$wsRD.Range("B10").Value = 20.69
$wsRD.Range("A11").Value = "TN"
$wsRD.Range("B11").Value = 98.58
$wsRD.Range("A12").Value = "JN"
$wsRD.Range("B12").Value = 47.22
$wsRD.Range("B13").Value = -3.59
$wsRD.Range("B14").Value = 12.45

# ---------------------------------------------------------------------------
# Sheet: Platform_Summary  (KPI update for platforms / total)
# ---------------------------------------------------------------------------
$wsPlatform = $wb.Worksheets.Item("Platform_Summary")
$wsPlatform.Range("B2").Value = 1442.8799999999999
$wsPlatform.Range("B5").Value = 3679.42
$wsPlatform.Range("B6").Value = 5227.44
